$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain decimal numbers need to be forced to text
# (so Excel keeps the exact original formatted string instead of converting to a float),
# by temporarily applying a text number format, then reverting the style afterwards.
$textCells = "D4","D5","D6","D7","D8","D9","D10","D11","D12","D15","D17","D19","D20","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D32","D33","D35","D36","D37","D38","D40","D41","D42","D43","D44","D48","D49","D50"
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '25.959.04'
$ws.Range("E2").Value = '  -1.45%  '

$ws.Range("D3").Value = '1.639.03'
$ws.Range("E3").Value = '  -1.62%  '

$ws.Range("D4").Value = '1.007'

$ws.Range("D5").Value = '215.09'
$ws.Range("E5").Value = '  -1.66%  '

$ws.Range("D6").Value = '0.5039'
$ws.Range("E6").Value = '  -2.38%  '

$ws.Range("D7").Value = '1.007'
$ws.Range("E7").Value = '  -0.13%  '

$ws.Range("B8").Value = 'Dogecoin'
$ws.Range("C8").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D8").Value = '0.06444'
$ws.Range("E8").Value = '  -0.10%  '

$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").Value = '0.2572'
$ws.Range("E9").Value = '  +0.21%  '

$ws.Range("D10").Value = '19.52'
$ws.Range("E10").Value = '  -1.93%  '

$ws.Range("D11").Value = '0.07737'
$ws.Range("E11").Value = '  +1.00%  '

$ws.Range("D12").Value = '4.253'
$ws.Range("E12").Value = '  -1.32%  '

$ws.Range("D13").Value = '1.631.22'
$ws.Range("E13").Value = '  -2.23%  '

$ws.Range("D14").Value = '1.863.91'
$ws.Range("E14").Value = '  -1.66%  '

$ws.Range("D15").Value = '0.5454'
$ws.Range("E15").Value = '  -1.66%  '

$ws.Range("D16").Value = '0.0₅7951'
$ws.Range("E16").Value = '  -1.10%  '

$ws.Range("D17").Value = '63.55'
$ws.Range("E17").Value = '  -1.43%  '

$ws.Range("D18").Value = '25.966.44'
$ws.Range("E18").Value = '  -1.54%  '

$ws.Range("D19").Value = '1.007'
$ws.Range("E19").Value = '  -0.14%  '

$ws.Range("D20").Value = '204.58'
$ws.Range("E20").Value = '  -2.74%  '

$ws.Range("D21").Value = '4.311'
$ws.Range("E21").Value = '  -2.05%  '

$ws.Range("D22").Value = '9.988'
$ws.Range("E22").Value = '  -1.19%  '

$ws.Range("D23").Value = '5.957'
$ws.Range("E23").Value = '  +0.95%  '

$ws.Range("D24").Value = '1.008'

$ws.Range("D25").Value = '1.906'
$ws.Range("E25").Value = '  +8.62%  '

$ws.Range("D26").Value = '141.18'
$ws.Range("E26").Value = '  -2.32%  '

$ws.Range("D27").Value = '0.1155'
$ws.Range("E27").Value = '  -0.36%  '

$ws.Range("D28").Value = '15.79'
$ws.Range("E28").Value = '  +0.18%  '

$ws.Range("D29").Value = '6.747'
$ws.Range("E29").Value = '  -3.48%  '

$ws.Range("D30").Value = '0.05060'
$ws.Range("E30").Value = '  -3.70%  '

$ws.Range("E31").Value = '  -1.85%  '

$ws.Range("D32").Value = '3.265'
$ws.Range("E32").Value = '  -3.31%  '

$ws.Range("D33").Value = '3.193'
$ws.Range("E33").Value = '  -0.97%  '

$ws.Range("E34").Value = '  -1.56%  '

$ws.Range("D35").Value = '2.338'
$ws.Range("E35").Value = '  -1.74%  '

$ws.Range("D36").Value = '0.8944'
$ws.Range("E36").Value = '  -3.32%  '

$ws.Range("D37").Value = '2.620'
$ws.Range("E37").Value = '  -5.04%  '

$ws.Range("D38").Value = '0.5656'
$ws.Range("E38").Value = '  -1.26%  '

$ws.Range("D39").Value = '1.144.14'
$ws.Range("E39").Value = '  -1.10%  '

$ws.Range("D40").Value = '0.01566'
$ws.Range("E40").Value = '  -1.99%  '

$ws.Range("D41").Value = '2.563'
$ws.Range("E41").Value = '  -0.63%  '

$ws.Range("D42").Value = '1.007'
$ws.Range("E42").Value = '  -0.19%  '

$ws.Range("D43").Value = '5.636'
$ws.Range("E43").Value = '  -0.30%  '

$ws.Range("D44").Value = '0.8183'
$ws.Range("E44").Value = '  -2.94%  '

$ws.Range("E45").Value = '  -0.49%  '

$ws.Range("D46").Value = '1.777.95'
$ws.Range("E46").Value = '  -1.53%  '

$ws.Range("D47").Value = '0.0₈112'
$ws.Range("E47").Value = '  -1.90%  '

$ws.Range("D48").Value = '0.4523'
$ws.Range("E48").Value = '  +0.41%  '

$ws.Range("D49").Value = '1.009'
$ws.Range("E49").Value = '  +0.22%  '

$ws.Range("D50").Value = '54.79'
$ws.Range("E50").Value = '  -2.36%  '

$ws.Range("E51").Value = '  -1.38%  '

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}